# Update the "dSF" column (F) values for a set of rows to re-pulled /
# re-pushed data per the "repull data, push all data, mean calculation"
# commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = 1
$ws.Range("F6").Value = 2
$ws.Range("F13").Value = 2
$ws.Range("F15").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("F24").Value = 1
$ws.Range("F27").Value = 1
$ws.Range("F32").Value = -1
$ws.Range("F33").Value = -1
$ws.Range("F36").Value = -13
$ws.Range("F41").Value = 2
$ws.Range("F42").Value = -1
$ws.Range("F59").Value = 0
$ws.Range("F60").Value = -6
$ws.Range("F61").Value = -1
$ws.Range("F64").Value = -1
$ws.Range("F72").Value = 4
$ws.Range("F73").Value = -1
$ws.Range("F76").Value = -13
